# Applies the "Implemented Sep 29 feedback" edit to todi_instructions.docx
$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $ok = $d.Content.Find.Execute($find, $true, $true, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Host "WARNING: find failed for: $find"
    }
    return $ok
}

# 1. Merge "How to " + "use a Transfer on Death Instrument" run split (no visible text change,
#    but normalize anyway via a no-op style replace so downstream finds are unaffected).
Replace-Text "How to use a Transfer on Death Instrument" "How to use a Transfer on Death Instrument"

# 2. Add the word "both" so notary sentence reads "...confirm that you both signed under oath."
Replace-Text "confirm that you signed under oath." "confirm that you both signed under oath."

# 3. Add missing period after "...document number, if known"
Replace-Text "the document number, if known" "the document number, if known."

# 4. Rewrite the "After your death..." paragraph to add the after_both_death merge field and
#    split the sentence about signing/recording.
Replace-Text "After your death, the beneficiaries must complete the Notice of Death Affidavit and Acceptance of Transfer on Death Instrument and then take the forms to the same Recorder of Deeds office where the TODI was recorded. This must be done within 2 years of " "After {{after_both_death}}, the beneficiaries must complete the Notice of Death Affidavit and Acceptance of Transfer on Death Instrument . They must sign it before 2 witnesses and a notary public. Then they must take the Notice to the same Recorder of Deeds office where the TODI was recorded.This must be done within 2 years of "

# 5. Give the trailing empty paragraph (after the table) explicit Arial paragraph-mark formatting.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.Font.NameAscii = "Arial"
$lastPara.Range.Font.NameOther = "Arial"
$lastPara.Range.Font.NameBi = "Arial"
